$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the leading index column (A), shifting Category/Sales/contribution % left.
$ws.Range("A1").EntireColumn.Delete()
